# Applies the "add ti_admin technical user + new user row" edit described
# by the commit: a new admin user login is introduced (H4: kcadmin -> ti_admin)
# and a brand new user row (row 5) is appended to the "utilisateurs" sheet,
# complete with a mailto hyperlink on the e-mail cell, mirroring the layout
# already used by existing rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("utilisateurs")

# --- H4: the "system" admin login changes from kcadmin to ti_admin ---------
$ws.Range("H4").Value = "ti_admin"

# --- Row 5: brand new user --------------------------------------------------
$ws.Range("A5").Value = "B"
$ws.Range("B5").Value = "ANOTHER"
$ws.Range("C5").Value = "77"
$ws.Range("D5").Value = "Un mec pas de l" + [char]0x2019 + "URSSAF"
$ws.Range("E5").Value = "NO"
$ws.Range("F5").Value = "quelqun"
$ws.Range("G5").Value = "pasdelurssaf"
$ws.Range("L5").Value = "pasdelurssaf.fr"
$ws.Range("M5").Value = 1

# Match the formatting already used for column C/D (text number format) and
# the rest of the row (default format), by copying from the row above.
$ws.Range("C2:D2").Copy()
$ws.Range("C5:D5").PasteSpecial(-4122)

$ws.Range("A2:B2").Copy()
$ws.Range("A5:B5").PasteSpecial(-4122)

$ws.Range("E2:G2").Copy()
$ws.Range("E5:G5").PasteSpecial(-4122)

$ws.Range("M2").Copy()
$ws.Range("M5").PasteSpecial(-4122)

# Email cell + hyperlink, styled like the other mail cells (H2/H4).
$ws.Hyperlinks.Add($ws.Range("H5"), "mailto:quelqun@pasdelurssaf.fr", $null, $null, "quelqun@pasdelurssaf.fr")
$ws.Range("H2").Copy()
$ws.Range("H5").PasteSpecial(-4122)

# Column L (DOMAINE MAIL, now used by the new row) gets an explicit width,
# and column C keeps its width but becomes an explicit (custom) one, exactly
# like a normal "I typed in column L/C" side effect in Excel.
$ws.Columns.Item(12).ColumnWidth = 13.69
$ws.Columns.Item(3).ColumnWidth = 10.69

# Leave the cursor on H4, matching where the edit session ended up.
$ws.Range("H4").Select() | Out-Null

Write-Output "edit applied"
